$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.216143727302551
$ws.Range("B1").Value = 2.281803369522095
$ws.Range("C1").Value = 3.457645416259766
$ws.Range("D1").Value = 2.181242227554321
$ws.Range("E1").Value = 1.320962071418762
